# Generate Report for Handback
# The file "4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md" has now been handed
# back successfully (in sync with en-US) for both zh-cn and de-de, so
# update the status / datetime / error-detail columns across all three
# sheets to reflect the new handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the 4f38d11f-...md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the 4f38d11f-...md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-15 14:46:31"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the 4f38d11f-...md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-15 14:46:39"
$wsDeDe.Range("P3").Value = ""
